$d = $word.ActiveDocument

# The document is a numbered/bulleted (numId=2) list of single-letter
# test paragraphs: A / " B" / C / " D" / "A " / B / "C " / D / A / B.
# This edit adds a handful of "dirty" whitespace / empty-line test
# cases among them, per the commit message:
#   "test: include test cases for trailing whitespace, empty lines etc"
#
# Net effect (1-based Word paragraph numbers, unchanged paragraphs
# omitted):
#   3: "C"  -> ""     (now a wholly empty list item / empty line)
#   5: "A " -> " "    (drop the leading letter, keep the trailing space)
#   6: "B"  -> ""     (now a wholly empty list item / empty line)
#   8: "D"  -> "C"    (swap the letter, keeping its paragraph spacing)

# --- Paragraph 3: "C" -> "" ------------------------------------------
# A plain Range.Text = "" leaves a placeholder <w:t/> behind; deleting
# the whole paragraph (through its mark) and then inserting a fresh
# empty paragraph in its place reproduces a genuinely empty line that
# still inherits the surrounding list/style formatting.
$p = $d.Paragraphs(3)
$p.Range.Delete()
$p = $d.Paragraphs(3)
$p.Range.InsertParagraphBefore()

# --- Paragraph 5: "A " -> " " (trailing space kept) -------------------
$p = $d.Paragraphs(5)
$p.Range.Text = " "

# --- Paragraph 6: "B" -> "" --------------------------------------------
$p = $d.Paragraphs(6)
$p.Range.Delete()
$p = $d.Paragraphs(6)
$p.Range.InsertParagraphBefore()

# --- Paragraph 8: "D" -> "C" --------------------------------------------
$p = $d.Paragraphs(8)
$p.Range.Text = "C"
